# ------------------------------------------------------------------------
# Adds a calculated "insert" column (E) to the Tabela1 table that builds an
# INSERT statement from the "NLQ preproccessed by GLAMORISE" and
# "NLIDB SQL" columns, and appends four new data rows (19-22) to the table.
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- 1. Add the new calculated column "insert" to the table -------------
$newCol = $tbl.ListColumns.Add()
$ws.Range("E1").Value = "insert"

$formula = '="INSERT INTO NLIDB_SQL_FROM_NLQ  VALUES(''"&Tabela1[[#This Row],[NLQ preproccessed by GLAMORISE]]&"'', ''"&SUBSTITUTE(Tabela1[[#This Row],[NLIDB SQL]],"''","''''")&"'');"'

# Fill the formula one row at a time (rows 2-18) so every cell keeps its own
# full formula text instead of being collapsed into a single shared formula.
for ($r = 2; $r -le 18; $r++) {
    $ws.Range("E$r").Formula = $formula
}

# --- 2. Append the four new rows of NLQ/SQL data -------------------------
$rowsData = @()
$rowsData += ,@('Which field produces the most oil per month?', 'Which field produces the oil month?', 'SELECT field, oil_production, year, month FROM ANP ', 'SELECT year, month, field, max(oil_production) as max_oil_production FROM NLIDB_RESULT_SET GROUP BY year, month, field ORDER BY year, month, field')
$rowsData += ,@('Which basin has the highest yearly oil production?', 'Which basin has the year oil production?', 'SELECT basin, year, oil_production FROM ANP ', 'SELECT basin, max(oil_production) as max_oil_production FROM (SELECT basin, SUM(oil_production) as oil_production FROM NLIDB_RESULT_SET GROUP BY basin, year) GROUP BY basin ORDER BY basin')
$rowsData += ,@('Which federated state has the lowest gas production?', 'Which federated state has the gas production?', 'SELECT state, gas_production FROM ANP ', 'SELECT state, min(gas_production) as min_gas_production FROM NLIDB_RESULT_SET GROUP BY state ORDER BY state')
$rowsData += ,@('Which state of the federation has the lowest gas production?', 'Which state of the federation has the gas production?', 'SELECT state, gas_production FROM ANP ', 'SELECT state, min(gas_production) as min_gas_production FROM NLIDB_RESULT_SET GROUP BY state ORDER BY state')

$rowIndex = 19
foreach ($row in $rowsData) {
    $tbl.ListRows.Add() | Out-Null
    $ws.Range("A$rowIndex").Value = $row[0]
    $ws.Range("B$rowIndex").Value = $row[1]
    $ws.Range("C$rowIndex").Value = $row[2]
    $ws.Range("D$rowIndex").Value = $row[3]
    $ws.Range("E$rowIndex").Formula = $formula
    $rowIndex++
}

# --- 3. Restore view state: scroll to column E and select E19:E22 --------
$ws.Range("E19:E22").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
